# Auto commit at 2025-10-13  7:47:43.56
#
# Appends one day's worth of readings (2025-10-12, serial 45942) for both
# charging stations to the bottom of the log on Sheet1: row 84 for
# "四方坪站" (station code shared-string index 2) and row 85 for "高岭站"
# (shared-string index 3). Each new row mirrors the formula pattern used
# by every preceding day block (see rows 82/83 for the immediately
# preceding day): column C = total minutes / gun count, D = C in day
# units, E = total kWh / gun count, F = total revenue / gun count,
# G = kWh per hour, H = average charges per gun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 84: 四方坪站 (2025-10-12) ----
$ws.Range("A84").Value = 45942
$ws.Range("B84").Value = "四方坪站"
$ws.Range("C84").Formula = "=16008/126"
$ws.Range("D84").Formula = "=C84/(24*60)"
$ws.Range("E84").Formula = "=9225.03/126"
$ws.Range("F84").Formula = "=3175.58/126"
$ws.Range("G84").Formula = "=9225.03/(16008/60)"
$ws.Range("H84").Formula = "=375/126"

# ---- Row 85: 高岭站 (2025-10-12) ----
$ws.Range("A85").Value = 45942
$ws.Range("B85").Value = "高岭站"
$ws.Range("C85").Formula = "=5235/36"
$ws.Range("D85").Formula = "=C85/(24*60)"
$ws.Range("E85").Formula = "=3777.63/36"
$ws.Range("F85").Formula = "=991.63/36"
$ws.Range("G85").Formula = "=3777.63/(5235/60)"
$ws.Range("H85").Formula = "=134/36"

# ---- Scroll / selection bookkeeping, matching the author's view state ----
$win = $excel.ActiveWindow
$win.ScrollRow = 25      # topLeftCell row -> B25
$win.ScrollColumn = 2    # topLeftCell col -> B25
$ws.Range("G85").Select()
